$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 192, shifting existing rows 192-254 down to 193-255
$ws.Rows(192).Insert()

# Populate the newly inserted row 192 with the new price record
$ws.Cells.Item(192, 1).Value = 3
$ws.Cells.Item(192, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(192, 3).Value = "Coquimbo"
$ws.Cells.Item(192, 4).Value = 44588
$ws.Cells.Item(192, 5).Value = 5
$ws.Cells.Item(192, 6).Value = 100112039
$ws.Cells.Item(192, 7).Value = "Ciboulette"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 180
$ws.Cells.Item(192, 11).Value = 1500
$ws.Cells.Item(192, 12).Value = 1500
$ws.Cells.Item(192, 13).Value = 1500
$ws.Cells.Item(192, 14).Value = "$/docena de atados"
$ws.Cells.Item(192, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(192, 16).Value = 500
$ws.Cells.Item(192, 17).Value = 3
$ws.Cells.Item(192, 18).Value = "Hortaliza"
